# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
# Both sheets list the same events; "全部类型" has one extra row inserted
# before row 32 ("苏州·爱乐之城..."), which shifts the later rows down by
# one on that sheet relative to "展览".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Row -> new F value, as it appears on the "展览" sheet.
$updates1 = @{
    3  = 560
    4  = 1130
    6  = 71
    9  = 1159
    10 = 16293
    11 = 275
    14 = 6353
    15 = 637
    18 = 21
    20 = 1271
    21 = 31
    30 = 5048
    32 = 11305
    35 = 147
    36 = 203
    38 = 270
}

# Same updates, with rows as they appear on the "全部类型" sheet (shifted by
# +1 for rows after the extra entry at row 32 on that sheet).
$updates4 = @{
    3  = 560
    4  = 1130
    6  = 71
    9  = 1159
    10 = 16293
    11 = 275
    14 = 6353
    15 = 637
    18 = 21
    20 = 1271
    21 = 31
    30 = 5048
    33 = 11305
    36 = 147
    37 = 203
    39 = 270
}

foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
